$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix shared string text for A.J GREEN (remove comma)
$ws.Range("A12").Value = "A.J GREEN"

# Update cell values in B and C columns (rows 2-11)
$ws.Range("C2").Value = 36

$ws.Range("B3").Value = 28
$ws.Range("C3").Value = 31

$ws.Range("C4").Value = 7

$ws.Range("B5").Value = 15
$ws.Range("C5").Value = 13

$ws.Range("B6").Value = 16
$ws.Range("C6").Value = 13

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0

$ws.Range("B8").Value = 17

$ws.Range("B9").Value = 11
$ws.Range("C9").Value = 17

$ws.Range("C10").Value = 2

$ws.Range("B11").Value = 5

# Update selection to C13
$ws.Range("C13").Select()
